# Auto-generated edit script applying cryptos.xlsx diff (cryptos price/volume update)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Ensure numeric-looking text values (e.g. prices like '228.43') are written
# as literal text, matching the source data's inline-string cells, instead of
# being auto-converted to numbers by Excel's input parser.
$numericCells = @()
$numericCells += ,@(5,4)
$numericCells += ,@(6,4)
$numericCells += ,@(8,4)
$numericCells += ,@(10,4)
$numericCells += ,@(13,4)
$numericCells += ,@(15,4)
$numericCells += ,@(16,4)
$numericCells += ,@(18,4)
$numericCells += ,@(20,4)
$numericCells += ,@(21,4)
$numericCells += ,@(22,4)
$numericCells += ,@(24,4)
$numericCells += ,@(25,4)
$numericCells += ,@(26,4)
$numericCells += ,@(27,4)
$numericCells += ,@(31,4)
$numericCells += ,@(32,4)
$numericCells += ,@(33,4)
$numericCells += ,@(34,4)
$numericCells += ,@(35,4)
$numericCells += ,@(37,4)
$numericCells += ,@(39,4)
$numericCells += ,@(40,4)
$numericCells += ,@(42,4)
$numericCells += ,@(43,4)
$numericCells += ,@(44,4)
$numericCells += ,@(45,4)
$numericCells += ,@(47,4)
$numericCells += ,@(49,4)
$numericCells += ,@(50,4)
foreach ($coord in $numericCells) {
    $ws.Cells.Item($coord[0], $coord[1]).NumberFormat = '@'
}

$ws.Cells.Item(2, 4).Value = '35.514.69'
$ws.Cells.Item(2, 5).Value = '  +2.71%  '
$ws.Cells.Item(3, 4).Value = '1.849.24'
$ws.Cells.Item(3, 5).Value = '  +1.27%  '
$ws.Cells.Item(4, 5).Value = '  +0.26%  '
$ws.Cells.Item(5, 4).Value = '228.43'
$ws.Cells.Item(5, 5).Value = '  +0.68%  '
$ws.Cells.Item(6, 4).Value = '0.610'
$ws.Cells.Item(6, 5).Value = '  +2.24%  '
$ws.Cells.Item(7, 5).Value = '  +0.19%  '
$ws.Cells.Item(8, 4).Value = '41.91'
$ws.Cells.Item(8, 5).Value = '  +9.07%  '
$ws.Cells.Item(9, 5).Value = '  +4.92%  '
$ws.Cells.Item(10, 4).Value = '0.0692'
$ws.Cells.Item(10, 5).Value = '  +1.15%  '
$ws.Cells.Item(11, 5).Value = '  +2.78%  '
$ws.Cells.Item(12, 4).Value = '2.116.56'
$ws.Cells.Item(12, 5).Value = '  +1.30%  '
$ws.Cells.Item(13, 4).Value = '11.50'
$ws.Cells.Item(13, 5).Value = '  +1.12%  '
$ws.Cells.Item(14, 4).Value = '1.850.06'
$ws.Cells.Item(14, 5).Value = '  +1.12%  '
$ws.Cells.Item(15, 2).Value = 'Polkadot'
$ws.Cells.Item(15, 3).Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Cells.Item(15, 4).Value = '4.73'
$ws.Cells.Item(15, 5).Value = '  +5.51%  '
$ws.Cells.Item(16, 2).Value = 'Polygon'
$ws.Cells.Item(16, 3).Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Cells.Item(16, 4).Value = '0.668'
$ws.Cells.Item(16, 5).Value = '  +4.34%  '
$ws.Cells.Item(17, 4).Value = '35.507.79'
$ws.Cells.Item(17, 5).Value = '  +2.64%  '
$ws.Cells.Item(18, 4).Value = '70.06'
$ws.Cells.Item(18, 5).Value = '  +1.45%  '
$ws.Cells.Item(19, 2).Value = 'ShibaInu'
$ws.Cells.Item(19, 3).Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Cells.Item(19, 4).Value = '0.0₃0800'
$ws.Cells.Item(19, 5).Value = '  +2.32%  '
$ws.Cells.Item(20, 2).Value = 'BitcoinCash'
$ws.Cells.Item(20, 3).Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Cells.Item(20, 4).Value = '245.52'
$ws.Cells.Item(20, 5).Value = '  +0.29%  '
$ws.Cells.Item(21, 4).Value = '12.25'
$ws.Cells.Item(21, 5).Value = '  +7.93%  '
$ws.Cells.Item(22, 4).Value = '4.78'
$ws.Cells.Item(22, 5).Value = '  +14.54%  '
$ws.Cells.Item(23, 5).Value = '  +0.24%  '
$ws.Cells.Item(24, 4).Value = '2.22'
$ws.Cells.Item(24, 5).Value = '  -0.59%  '
$ws.Cells.Item(25, 4).Value = '172.14'
$ws.Cells.Item(25, 5).Value = '  +0.03%  '
$ws.Cells.Item(26, 4).Value = '7.86'
$ws.Cells.Item(26, 5).Value = '  -1.33%  '
$ws.Cells.Item(27, 4).Value = '17.78'
$ws.Cells.Item(27, 5).Value = '  -0.46%  '
$ws.Cells.Item(28, 5).Value = '  +1.73%  '
$ws.Cells.Item(29, 5).Value = '  +0.22%  '
$ws.Cells.Item(30, 4).Value = '3.411.37'
$ws.Cells.Item(30, 5).Value = '  +40.40%  '
$ws.Cells.Item(31, 4).Value = '1.36'
$ws.Cells.Item(31, 5).Value = '  +10.18%  '
$ws.Cells.Item(32, 2).Value = 'Filecoin'
$ws.Cells.Item(32, 3).Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Cells.Item(32, 4).Value = '3.92'
$ws.Cells.Item(32, 5).Value = '  +2.18%  '
$ws.Cells.Item(33, 2).Value = 'InternetComputer(DFINITY)'
$ws.Cells.Item(33, 3).Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Cells.Item(33, 4).Value = '4.04'
$ws.Cells.Item(33, 5).Value = '  +2.08%  '
$ws.Cells.Item(34, 4).Value = '0.0536'
$ws.Cells.Item(34, 5).Value = '  +1.97%  '
$ws.Cells.Item(35, 4).Value = '1.87'
$ws.Cells.Item(35, 5).Value = '  +1.21%  '
$ws.Cells.Item(36, 5).Value = '  +3.29%  '
$ws.Cells.Item(37, 2).Value = 'TrustWalletToken'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Cells.Item(37, 4).Value = '1.09'
$ws.Cells.Item(37, 5).Value = '  +1.64%  '
$ws.Cells.Item(38, 2).Value = 'Maker'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Cells.Item(38, 4).Value = '1.340.93'
$ws.Cells.Item(38, 5).Value = '  -2.32%  '
$ws.Cells.Item(39, 2).Value = 'ARBITRUM'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Cells.Item(39, 4).Value = '1.03'
$ws.Cells.Item(39, 5).Value = '  +6.65%  '
$ws.Cells.Item(40, 2).Value = 'Aave'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Cells.Item(40, 4).Value = '88.14'
$ws.Cells.Item(40, 5).Value = '  +7.50%  '
$ws.Cells.Item(41, 5).Value = '  +2.93%  '
$ws.Cells.Item(42, 4).Value = '2.43'
$ws.Cells.Item(42, 5).Value = '  +1.06%  '
$ws.Cells.Item(43, 4).Value = '1.29'
$ws.Cells.Item(43, 5).Value = '  +6.13%  '
$ws.Cells.Item(44, 4).Value = '15.13'
$ws.Cells.Item(44, 5).Value = '  +7.07%  '
$ws.Cells.Item(45, 4).Value = '2.46'
$ws.Cells.Item(45, 5).Value = '  +1.02%  '
$ws.Cells.Item(46, 5).Value = '  +0.72%  '
$ws.Cells.Item(47, 4).Value = '0.0522'
$ws.Cells.Item(47, 5).Value = '  +3.14%  '
$ws.Cells.Item(48, 2).Value = 'RocketPoolETH'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Cells.Item(48, 4).Value = '2.016.97'
$ws.Cells.Item(48, 5).Value = '  +1.36%  '
$ws.Cells.Item(49, 2).Value = 'FraxShare'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Cells.Item(49, 4).Value = '6.03'
$ws.Cells.Item(49, 5).Value = '  +3.13%  '
$ws.Cells.Item(50, 4).Value = '104.50'
$ws.Cells.Item(50, 5).Value = '  +1.12%  '
$ws.Cells.Item(51, 5).Value = '  +0.19%  '
